$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Total row: only the numeric totals change (B2 stays blank)
$ws.Range("D2").Value = 14
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 29

# Row 3
$ws.Range("B3").Value = "Graduate or professional degree (MA, MS, MBA, PhD, JD, MD, DDS etc.)"
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 24
$ws.Range("G3").Value = 82.8

# Row 4
$ws.Range("B4").Value = "Bachelor’s degree"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 13.8

# Row 5
$ws.Range("B5").Value = "Other (please specify below)"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.4

# Row 6
$ws.Range("B6").Value = "Academia"
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 26
$ws.Range("G6").Value = 89.7

# Row 7 - only G changes
$ws.Range("G7").Value = 6.9

# Row 8
$ws.Range("B8").Value = "Government"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.4

# Row 9
$ws.Range("B9").Value = "Postdoc"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = 13.8

# Row 10
$ws.Range("B10").Value = "Other (please specify below)"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 6.9

# Row 11
$ws.Range("B11").Value = "Graduate student (including professional school student)"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 13.8

# Row 12
$ws.Range("B12").Value = "Undergraduate student"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.4

# Row 13
$ws.Range("B13").Value = "Staff member (including research/academic/teaching staff)"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 9
$ws.Range("G13").Value = 31

# Row 14
$ws.Range("B14").Value = "Faculty member"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 5
$ws.Range("F14").Value = 9
$ws.Range("G14").Value = 31

# Row 15
$ws.Range("B15").Value = "More than 5 years"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 9
$ws.Range("G15").Value = 31

# Row 16
$ws.Range("B16").Value = "1 to 5 years"
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 51.7

# Row 17 - only G changes
$ws.Range("G17").Value = 17.2
